$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '35.431.92'
$ws.Cells.Item(2, 5).Value = '  +0.54%  '

$ws.Cells.Item(3, 4).Value = '1.918.25'
$ws.Cells.Item(3, 5).Value = '  +1.26%  '

$ws.Cells.Item(4, 5).Value = '  +0.21%  '

$ws.Cells.Item(5, 4).Value = '''0.725'
$ws.Cells.Item(5, 5).Value = '  +10.65%  '

$ws.Cells.Item(6, 4).Value = '''253.53'
$ws.Cells.Item(6, 5).Value = '  +4.09%  '

$ws.Cells.Item(7, 5).Value = '  +0.23%  '

$ws.Cells.Item(8, 4).Value = '''40.71'
$ws.Cells.Item(8, 5).Value = '  -1.50%  '

$ws.Cells.Item(9, 5).Value = '  +3.03%  '

$ws.Cells.Item(10, 4).Value = '''52.48'
$ws.Cells.Item(10, 5).Value = '  +3.18%  '

$ws.Cells.Item(11, 4).Value = '''0.0748'
$ws.Cells.Item(11, 5).Value = '  +5.33%  '

$ws.Cells.Item(12, 4).Value = '''0.0992'
$ws.Cells.Item(12, 5).Value = '  -0.11%  '

$ws.Cells.Item(13, 4).Value = '2.198.89'
$ws.Cells.Item(13, 5).Value = '  +1.49%  '

$ws.Cells.Item(14, 4).Value = '''12.75'
$ws.Cells.Item(14, 5).Value = '  +6.69%  '

$ws.Cells.Item(15, 4).Value = '''0.723'
$ws.Cells.Item(15, 5).Value = '  +4.63%  '

$ws.Cells.Item(16, 4).Value = '1.922.24'
$ws.Cells.Item(16, 5).Value = '  +1.54%  '

$ws.Cells.Item(17, 4).Value = '''4.93'
$ws.Cells.Item(17, 5).Value = '  +1.92%  '

$ws.Cells.Item(18, 4).Value = '35.439.52'
$ws.Cells.Item(18, 5).Value = '  +0.60%  '

$ws.Cells.Item(19, 4).Value = '''74.38'
$ws.Cells.Item(19, 5).Value = '  +4.48%  '

$ws.Cells.Item(20, 4).Value = '0.0₃0843'
$ws.Cells.Item(20, 5).Value = '  +3.49%  '

$ws.Cells.Item(21, 2).Value = 'Avalanche'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Cells.Item(21, 4).Value = '''13.12'
$ws.Cells.Item(21, 5).Value = '  +5.73%  '

$ws.Cells.Item(22, 2).Value = 'BitcoinCash'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Cells.Item(22, 4).Value = '''243.17'
$ws.Cells.Item(22, 5).Value = '  +1.11%  '

$ws.Cells.Item(23, 4).Value = '''5.13'
$ws.Cells.Item(23, 5).Value = '  +7.92%  '

$ws.Cells.Item(24, 5).Value = '  +0.07%  '

$ws.Cells.Item(25, 4).Value = '''2.42'
$ws.Cells.Item(25, 5).Value = '  +0.42%  '

$ws.Cells.Item(26, 5).Value = '  +2.24%  '

$ws.Cells.Item(27, 4).Value = '''167.33'
$ws.Cells.Item(27, 5).Value = '  -1.72%  '

$ws.Cells.Item(28, 4).Value = '''8.67'
$ws.Cells.Item(28, 5).Value = '  +2.96%  '

$ws.Cells.Item(29, 5).Value = '  +6.09%  '

$ws.Cells.Item(30, 4).Value = '''18.75'
$ws.Cells.Item(30, 5).Value = '  +2.71%  '

$ws.Cells.Item(31, 4).Value = '4.126.10'

$ws.Cells.Item(32, 5).Value = '  +6.24%  '

$ws.Cells.Item(33, 4).Value = '''2.02'
$ws.Cells.Item(33, 5).Value = '  +16.34%  '

$ws.Cells.Item(34, 4).Value = '''1.64'
$ws.Cells.Item(34, 5).Value = '  +23.66%  '

$ws.Cells.Item(35, 4).Value = '''0.0582'
$ws.Cells.Item(35, 5).Value = '  +3.94%  '

$ws.Cells.Item(36, 5).Value = '  +3.06%  '

$ws.Cells.Item(37, 5).Value = '  +0.18%  '

$ws.Cells.Item(38, 5).Value = '  -3.81%  '

$ws.Cells.Item(39, 5).Value = '  +0.71%  '

$ws.Cells.Item(40, 4).Value = '''17.47'
$ws.Cells.Item(40, 5).Value = '  +8.16%  '

$ws.Cells.Item(41, 4).Value = '''96.86'
$ws.Cells.Item(41, 5).Value = '  +8.49%  '

$ws.Cells.Item(42, 5).Value = '  +3.34%  '

$ws.Cells.Item(43, 5).Value = '  +0.78%  '

$ws.Cells.Item(44, 4).Value = '''0.0657'
$ws.Cells.Item(44, 5).Value = '  +2.48%  '

$ws.Cells.Item(45, 4).Value = '1.343.50'
$ws.Cells.Item(45, 5).Value = '  +0.61%  '

$ws.Cells.Item(46, 4).Value = '''2.45'
$ws.Cells.Item(46, 5).Value = '  +2.83%  '

$ws.Cells.Item(47, 4).Value = '''2.42'
$ws.Cells.Item(47, 5).Value = '  +0.88%  '

$ws.Cells.Item(48, 4).Value = '''6.73'
$ws.Cells.Item(48, 5).Value = '  +2.90%  '

$ws.Cells.Item(49, 5).Value = '  +0.25%  '

$ws.Cells.Item(50, 5).Value = '  -6.23%  '

$ws.Cells.Item(51, 4).Value = '''11.74'
$ws.Cells.Item(51, 5).Value = '  +4.40%  '
